# Add remaining SOCP constraints
# - Fill in the previously-blank Qmax/Qmin (columns M and N) reactive-power
#   bound cells for the four VRE rows (8-11) with 0, completing the SOCP
#   (second-order cone) reactive-power constraint data alongside the
#   existing Qmax/Qmin=0 entries in column O.
# - Adjust the height of the two wrapped header rows (5 and 6) to fit.
# - Restore/update the active-cell selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns M (Qmax) and N (Qmin) for rows 8-11 were empty; set them to 0.
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0

$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0

$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0

$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0

# Resize the two wrapped-text header rows.
$ws.Rows.Item(5).RowHeight = 86
$ws.Rows.Item(6).RowHeight = 43

# Update the selected / active cell.
$null = $ws.Range("O18").Select()
